$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")

# Seminar of Oct. 16th (row 8): record attendance counts
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 6
